# Updates the cryptos list with the latest scraped price/volume data.
# Mirrors the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.345.39'
$ws.Range("D3").Value = '3.503.91'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Formula = "'585.83"
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Formula = "'135.76"
$ws.Range("E6").Value = '  +2.29%  '
$ws.Range("D7").Value = '3.504.59'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").Formula = "'0.125"
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("D11").Formula = "'7.16"
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("E12").Value = '  -3.70%  '
$ws.Range("D13").Value = '4.098.75'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '3.500.81'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '64.314.23'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Formula = "'25.54"
$ws.Range("E18").Value = '  -8.47%  '
$ws.Range("D19").Formula = "'9.77"
$ws.Range("E19").Value = '  -2.23%  '
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("E21").Value = '  -2.55%  '
$ws.Range("D22").Formula = "'383.57"
$ws.Range("E22").Value = '  -1.78%  '
$ws.Range("D23").Formula = "'0.570"
$ws.Range("E23").Value = '  -1.60%  '
$ws.Range("D24").Value = '3.639.85'
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").Formula = "'74.01"
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("D27").Formula = "'5.73"
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Formula = "'0.0000114"
$ws.Range("E28").Value = '  +3.64%  '
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("D31").Formula = "'0.999"
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("D34").Value = '3.521.06'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("E37").Value = '  -1.85%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  -2.29%  '
$ws.Range("D40").Formula = "'6.85"
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("D41").Formula = "'163.88"
$ws.Range("E41").Value = '  -4.32%  '
$ws.Range("D42").Formula = "'0.0785"
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("D44").Formula = "'25.84"
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Formula = "'4.41"
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").Formula = "'1.63"
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("D50").Value = '2.470.30'
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("D51").Formula = "'6.77"
$ws.Range("E51").Value = '  -2.06%  '

Write-Output "cryptos sheet refreshed: 81 cells updated"
